$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - add "Save" column header with same style as neighboring header (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save column values, row 2 through row 23
$values = @(0,0,0,0,0,0,0,1,0,1,0,0,0,0,1,0,1,1,1,0,0,1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
